$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("instrument_info")

# Add the new "instrument_function" worksheet after instrument_info
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "instrument_function"

# Header row
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "type"
$ws2.Range("C1").Value = "coeffs"

# Coefficients for each analog channel (A1CH1..A5CH8), row 2..41
$coeffs = @(
  "[1,1]",
  "[0,1]",
  "[2,1]",
  "[-1,1]",
  "[1,2]",
  "[0.5,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]",
  "[0,0.5]"
)

for ($i = 0; $i -lt 40; $i++) {
  $row = $i + 2
  $srcRow = $i + 82
  $name = $ws1.Cells.Item($srcRow, 1).Value()
  $ws2.Cells.Item($row, 1).Value = $name
  $ws2.Cells.Item($row, 2).Value = "poly"
  $ws2.Cells.Item($row, 3).Value = $coeffs[$i]
}

# Update selection on instrument_info (no longer the active tab)
$ws1.Range("H10").Select() | Out-Null

# Make instrument_function the active tab, with C8:C41 selected
$ws2.Select() | Out-Null
$ws2.Range("C8:C41").Select() | Out-Null
